$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D8").Value = "2016-03-10 05:00:15"
$wsZh.Range("G8").Value = "2016-03-10 05:01:18"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D8").Value = "2016-03-10 05:00:31"
$wsDe.Range("G8").Value = "2016-03-10 05:01:41"
